# edit.ps1 - apply the "create and setup login prefab and database scripts"
# changes to tsv_UI_Defaults.xlsx:
#   - insert 8 new rows (login/signup/email/password/game-key UI strings)
#     between the existing "btn main menu" row and the "UI wave" row
#   - re-format the whole table: vertically centered, word-wrapped, bordered
#   - tweak the sheet view (zoom 130%, selection on C11, scroll reset)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert 8 new rows right after row 5 ("UI button | btn main menu | ...")
#    Everything that used to be row 6 onward (UI wave, player stats, ...)
#    shifts down to row 14 onward automatically.
# ---------------------------------------------------------------------
$ws.Rows.Item(6).Resize(8).Insert()

# ---------------------------------------------------------------------
# 2. Fill in the new rows with the Login / Sign up / Email / Password /
#    Game key UI default strings.
# ---------------------------------------------------------------------
$newRows = @(
    @("UI button",         "btn validate",        "Validate",   "Valider"),
    @("UI Login",          "tmp login",            "Log in",     "Connexion"),
    @("UI Sign up",        "tmp signup",            "Sign up",    "Inscription"),
    @("UI email",          "tmp email",             "Email",      "Email"),
    @("UI password",       "tmp password",          "Password",   "Mot de passe"),
    @("UI confirm email",  "tmp confirm email",
        "To gain access to ZombieSurvivor's features, please verify your e-mail address.",
        "Pour accéder aux fonctionnalités de ZombieSurvivor, veuillez vérifier votre adresse e-mail."),
    @("UI Game key",       "tmp valid key",         "Valid key!",    "Clé valide!"),
    @("UI Game key",       "tmp invalid key",       "Invalid Key!",  "Clé invalide!")
)

$r = 6
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 3. Re-apply consistent formatting (thin black borders, vertical-center,
#    word-wrap) across the header row and the whole data table, including
#    the newly inserted rows.
# ---------------------------------------------------------------------
$fullRange = $ws.Range("A1:F26")
$fullRange.Borders.LineStyle = 1
$fullRange.Borders.Weight = 2
$fullRange.Borders.ColorIndex = 1
$fullRange.WrapText = $true
$fullRange.VerticalAlignment = -4108

# Row containing the long email-verification sentence needs extra height.
$ws.Rows.Item(11).RowHeight = 57.6

# ---------------------------------------------------------------------
# 4. Sheet view tweaks: zoom to 130%, reset scroll position, select C11.
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Application.ActiveWindow.Zoom = 130
$ws.Range("C11").Select()
